$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows to append starting at row 27 (data lineage at variable level,
# and inferred lineage at dataset level).
$rows = @(
    @{ Row=27; A=1747906147; B="update"; C="variable"; D="ser_pub_loc___variable_12"; F="source_var_ids"; H="ser_pub_loc___variable_10, ser_pub_loc___variable_11" },
    @{ Row=28; A=1747906281; B="add";    C="config";   D="alias_3" },
    @{ Row=29; A=1747906551; B="update"; C="variable"; D="ser_pub_loc___variable_12"; F="source_var_ids"; G="ser_pub_loc___variable_10, ser_pub_loc___variable_11" },
    @{ Row=30; A=1747906551; B="update"; C="variable"; D="ser_pub_loc___variable_12"; F="sourceVar_ids"; H="ser_pub_loc___variable_10, ser_pub_loc___variable_11" },
    @{ Row=31; A=1747906566; B="update"; C="config";   D="alias_3"; F="value"; G="variable : source_var"; H="variable : sourceVar" },
    @{ Row=32; A=1747907576; B="add";    C="config";   D="test" },
    @{ Row=33; A=1747907576; B="delete"; C="config";   D="alias_3" },
    @{ Row=34; A=1747907625; B="add";    C="config";   D="alias_3" },
    @{ Row=35; A=1747907625; B="delete"; C="config";   D="test" },
    @{ Row=36; A=1747913177; B="delete"; C="config";   D="alias_3" },
    @{ Row=37; A=1747913221; B="update"; C="variable"; D="ser_pub_loc___variable_12"; F="sourceVar_ids"; G="ser_pub_loc___variable_10, ser_pub_loc___variable_11" },
    @{ Row=38; A=1747913221; B="update"; C="variable"; D="ser_pub_loc___variable_12"; F="source_ids"; H="ser_pub_loc___variable_10, ser_pub_loc___variable_11" },
    @{ Row=39; A=1747913314; B="update"; C="variable"; D="ser_pub_loc___variable_12"; F="source_ids"; G="ser_pub_loc___variable_10, ser_pub_loc___variable_11" },
    @{ Row=40; A=1747913314; B="update"; C="variable"; D="ser_pub_loc___variable_12"; F="sourceVar_ids"; H="ser_pub_loc___variable_10, ser_pub_loc___variable_11" },
    @{ Row=41; A=1747918633; B="update"; C="variable"; D="ser_pub_loc___variable_13"; F="sourceVar_ids"; H="ser_pub_loc___variable_12" },
    @{ Row=42; A=1747921769; B="update"; C="variable"; D="ser_pub_loc___variable_14"; F="sourceVar_ids"; H="accident_route___variable_7, dep_sante___variable_7" },
    @{ Row=43; A=1747921769; B="update"; C="variable"; D="ser_pub_loc___variable_15"; F="sourceVar_ids"; H="accident_route___variable_7, dep_sante___variable_7" },
    @{ Row=44; A=1747921769; B="update"; C="variable"; D="dep_sante___variable_3";   F="sourceVar_ids"; H="accident_route___variable_7" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    if ($r.ContainsKey("B")) { $ws.Cells.Item($row, 2).Value = $r.B }
    if ($r.ContainsKey("C")) { $ws.Cells.Item($row, 3).Value = $r.C }
    if ($r.ContainsKey("D")) { $ws.Cells.Item($row, 4).Value = $r.D }
    if ($r.ContainsKey("E")) { $ws.Cells.Item($row, 5).Value = $r.E }
    if ($r.ContainsKey("F")) { $ws.Cells.Item($row, 6).Value = $r.F }
    if ($r.ContainsKey("G")) { $ws.Cells.Item($row, 7).Value = $r.G }
    if ($r.ContainsKey("H")) { $ws.Cells.Item($row, 8).Value = $r.H }
    if ($r.ContainsKey("I")) { $ws.Cells.Item($row, 9).Value = $r.I }
}
